$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cells / formulas -------------------------------------------------
# Row 7: powers-of-two helper columns
$ws.Range("I7").Formula = "=2^28"
$ws.Range("J7").Formula = "=I7-H7"

# Row 13: ratios against H7
$ws.Range("H13").Formula = "=H7/H7"
$ws.Range("I13").Formula = "=I7/H7"

# Row 30 / 36: further powers of two
$ws.Range("I30").Formula = "=2^29"
$ws.Range("I36").Formula = "=2^9"

# --- Number formatting: apply built-in "Comma [0]" style ------------------
# (this is what Excel does when you pick the Comma[0] style from the
#  Cell Styles gallery - it injects the numFmt/font/xf/cellStyle entries).
# Applied per-range because multi-area (union) ranges only style the first
# area in this host.
$ws.Range("F2:F25").Style = "Comma [0]"
$ws.Range("H7").Style = "Comma [0]"
$ws.Range("I7").Style = "Comma [0]"
$ws.Range("J7").Style = "Comma [0]"
$ws.Range("I30").Style = "Comma [0]"

# --- Column widths ----------------------------------------------------------
$ws.Columns("F").ColumnWidth = 12.166666666666666
$ws.Columns("H").ColumnWidth = 11.666666666666666
$ws.Columns("I").ColumnWidth = 11.833333333333334
$ws.Columns("J").ColumnWidth = 11.666666666666666

# --- Selection / view --------------------------------------------------------
$ws.Range("H19").Select()
